# Adapt column header formatting to respective input file names (FV2404 / FV2410)
# and add the Table1 ListObject + frozen header row, matching the commit's intent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells: "_old" -> "_FV2404", "_new" -> "_FV2410" -------
# Scope the Replace() to the header row only so no other cell content is touched.
$headerRange = $ws.Range("A1:U1")
$null = $headerRange.Replace("_old", "_FV2404", 2, 1, $false, $false, $false)
$null = $headerRange.Replace("_new", "_FV2410", 2, 1, $false, $false, $false)

# --- 2. Freeze the header row (pane split above row 2) ----------------------
$ws.Activate()
$null = $ws.Range("A1").Select()
$null = $ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$null = $ws.Range("A1").Select()

# --- 3. Turn the data range into an Excel Table named "Table1" -------------
$tableRange = $ws.Range("A1:U65")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
